$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New archival records to append below the existing header row (row 1).
# Columns: A=identifier, B=alternativeIdentifiers, C=title, D=date_s,
#          E=levelOfDescription, F=extentAndMedium, G=notes, H=file_path
$records = @(
    @{ id = "MCH158-1"; title = "KAIROS (HOLLAND), BOOJKS (DUTCH), PAMPHLETS (DUTCH)"; level = "Series"; extent = "1 Box"; notes = "LOCATION: 22A | GRAP COUNT NUMER: NONE" },
    @{ id = "MCH158-2"; title = "KAIROS (HOLLAND), BOOJKS (DUTCH), PAMPHLETS (DUTCH)"; level = "Series"; extent = "1 Box"; notes = "LOCATION: 22A | GRAP COUNT NUMER: NONE" },
    @{ id = "MCH158-3"; title = "BOOKS AND REPORTS, DOCUMENTS SERVICES";               level = "Series"; extent = "1 Box"; notes = "LOCATION: 22A | GRAP COUNT NUMER: NONE" },
    @{ id = "MCH158-4"; title = "BROCHURES, ETC";                                      level = "Series"; extent = "1 Box"; notes = "LOCATION: 22A | GRAP COUNT NUMER: NONE" }
)

$rowIndex = 2
foreach ($rec in $records) {
    # Apply the data-row font (Calibri 10, theme text color) to every used
    # column of this row except B (alternativeIdentifiers stays untouched/blank).
    $rowRange = $ws.Range("A$rowIndex" + ",C$rowIndex" + ":D$rowIndex" + ",E$rowIndex" + ":H$rowIndex")
    foreach ($area in $rowRange.Areas) {
        $area.Font.Name = "Calibri"
        $area.Font.Size = 10
        $area.Font.ThemeColor = 1
    }

    $ws.Range("A$rowIndex").Value = $rec.id
    $ws.Range("C$rowIndex").Value = $rec.title
    $ws.Range("E$rowIndex").Value = $rec.level
    $ws.Range("F$rowIndex").Value = $rec.extent
    $ws.Range("G$rowIndex").Value = $rec.notes

    $rowIndex = $rowIndex + 1
}

$lastRow = $rowIndex - 1

# Restore the frozen header pane (row 1 frozen) and select the newly added
# block of records, matching the worksheet view state of the saved file.
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$null = $ws.Range("A2:H$lastRow").Select()
